$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsToDelete = @(79, 75, 69, 52, 12, 11, 9, 8, 6, 4)
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}
